$d = $word.ActiveDocument

# Find the paragraph that contains the "Ok test" text we need to remove.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Ok test*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $targetPara = $d.Paragraphs.Item($targetIndex)
    $prevPara = $d.Paragraphs.Item($targetIndex - 1)

    # Remove the "Ok test" run text itself, but keep the paragraph mark and
    # whatever follows it (e.g. the _GoBack bookmark) intact.
    $textRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
    $textRange.Delete()

    # Now merge the (now empty) paragraph back into the previous one by
    # deleting the preceding paragraph mark. Anything that used to follow
    # the removed text (such as the bookmark) stays attached to the
    # surviving paragraph.
    $markRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
    $markRange.Delete()
}
